$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2432490129025711
$ws.Range("B3").Value = 0.2411897279234069
$ws.Range("B4").Value = 0.2279374252521792
$ws.Range("B5").Value = 0.2110671728913582
$ws.Range("B6").Value = 0.2232835019967376
$ws.Range("B7").Value = 0.2338366062972175
$ws.Range("B8").Value = 0.2135195243514688
$ws.Range("B9").Value = 0.200607347540649
$ws.Range("B10").Value = 0.2823875815070244
$ws.Range("B11").Value = 0.1998897476382224
$ws.Range("B12").Value = 0.2053661800737224
$ws.Range("B13").Value = 0.2025752555140319
$ws.Range("B14").Value = 0.1739384242826933
$ws.Range("B15").Value = 0.2015725473051372
$ws.Range("B16").Value = 0.2320889448780388
$ws.Range("B17").Value = 0.2167895105100233
$ws.Range("B18").Value = 0.2144442310725814
$ws.Range("B19").Value = 0.2121829174672534
$ws.Range("B20").Value = 0.1731373828639528
$ws.Range("B21").Value = 0.1940917271092404
$ws.Range("B22").Value = 0.2079469975294282
$ws.Range("B23").Value = 0.179211919645771
$ws.Range("B24").Value = 0.1674113558621942
$ws.Range("B25").Value = 0.1897901160543019
$ws.Range("B26").Value = 0.205923835628804
$ws.Range("B27").Value = 0.2205501952095101
$ws.Range("B28").Value = 0.2019318048492107
$ws.Range("B29").Value = 0.1958749450597683
$ws.Range("B30").Value = 0.1981261388742524
$ws.Range("B31").Value = 0.2391133932550645
$ws.Range("B32").Value = 0.1950598397699665
$ws.Range("B33").Value = 0.1487496675782086
$ws.Range("B34").Value = 0.1892665454765523
$ws.Range("B35").Value = 0.1765637362566932
$ws.Range("B36").Value = 0.2337421648130379
$ws.Range("B37").Value = 0.196159914135432
$ws.Range("B38").Value = 0.2096921815164276
$ws.Range("B39").Value = 0.1988977444286743
$ws.Range("B40").Value = 0.2290684984032781
$ws.Range("B41").Value = 0.2153576239808776
$ws.Range("B42").Value = 0.2127302161003796
$ws.Range("B43").Value = 0.2123330059331108
$ws.Range("B44").Value = 0.2111907239136736
$ws.Range("B45").Value = 0.2413968397084669
$ws.Range("B46").Value = 0.1838156528288608
$ws.Range("B47").Value = 0.2464650140770868
$ws.Range("B48").Value = 0.2234666136291428
$ws.Range("B49").Value = 0.1959046171534852
$ws.Range("B50").Value = 0.2027371075584232
$ws.Range("B51").Value = 0.2207704047530455
$ws.Range("B52").Value = 0.2514823645950497
$ws.Range("B53").Value = 0.2046301203171215
$ws.Range("B54").Value = 0.1956943611878062
$ws.Range("B55").Value = 0.2466881491246491
$ws.Range("B56").Value = 0.2123710190639136
$ws.Range("B57").Value = 0.1820696510099117
$ws.Range("B58").Value = 0.210302979833812
$ws.Range("B59").Value = 0.1794839063464464
$ws.Range("B60").Value = 0.2233297240657419
$ws.Range("B61").Value = 0.2393133728178417
$ws.Range("B62").Value = 0.1629844592048781
$ws.Range("B63").Value = 0.1969827144392482
$ws.Range("B64").Value = 0.1965733061178754
$ws.Range("B65").Value = 0.183836812759245
$ws.Range("B66").Value = 0.1626435010973513
$ws.Range("B67").Value = 0.1995352458803386
$ws.Range("B68").Value = 0.2485047676408502
$ws.Range("B69").Value = 0.2187800333905454
$ws.Range("B70").Value = 0.2045486734926522
$ws.Range("B71").Value = 0.1906019554497276
$ws.Range("B72").Value = 0.2207784478459045
$ws.Range("B73").Value = 0.2263140766664
$ws.Range("B74").Value = 0.209776060412968
$ws.Range("B75").Value = 0.2082701329743925
$ws.Range("B76").Value = 0.2040152875554088
$ws.Range("B77").Value = 0.2138425143858722
$ws.Range("B78").Value = 0.1728289216916878
